# Updates cryptos list (price/volume refresh) as produced by the
# "Updated cryptos list ... with GitHub Actions" scraper run.
# For D-column cells whose new value is a plain numeric-looking string
# (e.g. "309.54", "0.4684", "1.670"), we force Text number format before
# assigning the value so Excel keeps it as a literal string (preserving
# trailing zeros / exact formatting) instead of silently converting it
# to a floating point number. The style is reset back to "Normal"
# afterwards so no stray cell formatting is introduced.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.893.49"
$ws.Range("E2").Value = "  +0.86%  "
$ws.Range("D3").Value = "1.818.00"
$ws.Range("E3").Value = "  +1.60%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "309.54"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.55%  "
$ws.Range("E6").Value = "  +0.05%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4684"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.95%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3690"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07377"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.63%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8710"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.05%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.43"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.28%  "
$ws.Range("D12").Value = "1.802.18"
$ws.Range("E12").Value = "  +0.64%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.364"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.72%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "92.47"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.49%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.07071"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.497"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.48%  "
$ws.Range("E17").Value = "  +0.03%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008703"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.34%  "
$ws.Range("E19").Value = "  -0.03%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.76"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.38%  "
$ws.Range("D21").Value = "26.949.11"
$ws.Range("E21").Value = "  +1.04%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.345"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.58%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.58"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.25%  "
$ws.Range("D24").Value = "2.071.50"
$ws.Range("E24").Value = "  +3.14%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.902"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.06%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "151.44"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.23%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.172"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.81%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.36"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.80%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.331"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.63%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "115.71"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.75%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08939"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.63%  "
$ws.Range("E32").Value = "  +2.28%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.164"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.92%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.511"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.88%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.903"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.79%  "
$ws.Range("E36").Value = "  +0.04%  "
$ws.Range("E37").Value = "  -1.60%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01966"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.74%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05283"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.73%  "
$ws.Range("B40").Value = "MXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.946"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.23%  "
$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.279"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.25%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5326"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.74%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.343"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.91%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1668"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.88%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.432"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.30%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4941"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.11%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.47"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.59%  "
$ws.Range("B48").Value = "PaxDollar"
$ws.Range("C48").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.001"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.10%  "
$ws.Range("B49").Value = "Quant"
$ws.Range("C49").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "103.91"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.21%  "
$ws.Range("B50").Value = "NEARProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.670"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.95%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06277"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.15%  "
